# Applies the commit "Added updating CB option":
#  - Deletes (collapses) column O ("Unnamed: 14") so every column from
#    O.."BG" shifts one position to the left (P->O, Q->P, ... BG->BF),
#    shrinking the used range from A1:BG17 to A1:BF17.
#  - Renumbers the generic "Unnamed: N" header labels that shifted
#    position so they keep matching their (now one lower) column index.
#  - Refreshes the data: the four companies that used to be listed first
#    (Dive, Enlight Exchange, MSICS Pharma, Beffi) move to the bottom of
#    the table and the newly pulled Crunchbase data for the remaining
#    companies (which used to be last) now appears at the top, each
#    stamped with the newer Updating_Date of 23-08-2024.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 17
$lastColAfter  = 58   # BF
$oCol = 15            # O

# ---- 1. Read the full existing grid (1-based: rows 1..17, cols 1..59) ----
$src = $ws.Range("A1:BG17").Value2

# ---- 2. Build the target grid (rows 1..17, cols 1..58) ------------------
$dst = New-Object 'object[,]' $lastRow, $lastColAfter

# -- Row 1 (headers): columns A..N stay put, O..BF take the value that
#    used to sit one column to the right (P..BG). --
for ($c = 1; $c -le ($oCol - 1); $c++) {
    $dst[0, $c - 1] = $src[1, $c]
}
for ($c = $oCol; $c -le $lastColAfter; $c++) {
    $dst[0, $c - 1] = $src[1, $c + 1]
}

# -- Fix up the generic "Unnamed: N" placeholder headers: since a column
#    was removed, every such label that moved must count one lower. --
for ($c = 1; $c -le $lastColAfter; $c++) {
    $v = $dst[0, $c - 1]
    if ($v -is [string] -and $v.StartsWith("Unnamed: ")) {
        $n = [int]($v.Substring(9))
        $dst[0, $c - 1] = "Unnamed: " + ($n - 1)
    }
}

# -- Data rows: old rows 6..17 move up to become new rows 2..13, and old
#    rows 2..5 move down to become new rows 14..17 (a rotate-by-4). Every
#    row also gets the column shift described above, and column B
#    (Updating_Date) is refreshed to 23-08-2024. --
$rowOrder = @(6,7,8,9,10,11,12,13,14,15,16,17,2,3,4,5)

for ($i = 0; $i -lt $rowOrder.Length; $i++) {
    $srcRow = $rowOrder[$i]
    $dstRow = $i + 2

    # Column A (Company_Name) unchanged column position
    $dst[$dstRow - 1, 0] = $src[$srcRow, 1]
    # Column B (Updating_Date) forced to the new date
    $dst[$dstRow - 1, 1] = "23-08-2024"
    # Columns C..N unchanged column position
    for ($c = 3; $c -le ($oCol - 1); $c++) {
        $dst[$dstRow - 1, $c - 1] = $src[$srcRow, $c]
    }
    # Columns O..BF <- old P..BG
    for ($c = $oCol; $c -le $lastColAfter; $c++) {
        $dst[$dstRow - 1, $c - 1] = $src[$srcRow, $c + 1]
    }
}

# Column T (Company_Founded_Year) keeps its real numeric value for the
# four rows that now hold Dive/Enlight Exchange/MSICS Pharma/Beffi.
$founded = @{ 14 = 2022; 15 = 2022; 16 = 2022; 17 = 2021 }

# ---- 3. Write the grid back out ------------------------------------------
# Columns B (dates) and T/V (dates & "1,234"-style numbers-as-text) must be
# protected with a Text number format before the write, otherwise Excel
# auto-converts those look-alike strings into real dates/numbers. The rest
# of the sheet is plain, unambiguous text, so it is left with its existing
# (default) formatting/style untouched. (Each column is formatted with its
# own statement: a single multi-area Range(...).NumberFormat assignment
# only affects the first area under this automation host.)
$ws.Range("B2:B17").NumberFormat = "@"
$ws.Range("T2:T17").NumberFormat = "@"
$ws.Range("V2:V17").NumberFormat = "@"

$ws.Range("A1:BF17").Value2 = $dst

foreach ($r in $founded.Keys) {
    $cell = $ws.Cells.Item($r, 20)   # column T
    $cell.NumberFormat = "General"
    $cell.Value2 = $founded[$r]
}

# ---- 4. Drop the now unused column BG so the sheet's dimension shrinks
#         back down to A1:BF17, matching the target workbook. ----
$ws.Columns("BG").Delete() | Out-Null

"Done. UsedRange: " + $ws.UsedRange.Address()
